# Commit: "add updated 2020 data"
#
# 1) "Population 2019-Corrected": Arizona (row 5) previously had no data
#    (blank cells, "not highlighted" style). Fill in the corrected AZ
#    population-admission figures and restyle the row to match the other
#    "updated/corrected" rows (fill-highlighted styles 53/54/56).
#
# 2) "Population 2020-Corrected": Arizona didn't have its own row yet in
#    this corrected sheet (only a hidden placeholder further down). Insert
#    a brand-new highlighted row for AZ right after Alaska/Alabama (before
#    the existing Arkansas row), populated with the 2020 corrected figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Population 2019-Corrected
# ---------------------------------------------------------------------
$popCorrected2020 = $wb.Worksheets.Item("Population 2020-Corrected")
$pop2019Corrected = $wb.Worksheets.Item("Population 2019-Corrected")

# Borrow the "highlighted / updated" formatting from a row that already
# carries it (row 2, Alaska) in the 2020-Corrected sheet, then paste just
# the formats onto the Arizona row (row 5) of the 2019-Corrected sheet.
$popCorrected2020.Range("A2:J2").Copy()
$pop2019Corrected.Range("A5:J5").PasteSpecial(-4122)
$popCorrected2020.Range("K2:M2").Copy()
$pop2019Corrected.Range("K5:M5").PasteSpecial(-4122)

$pop2019Corrected.Range("C5").Value = 41937
$pop2019Corrected.Range("D5").Value = 10080
$pop2019Corrected.Range("E5").Value = 8339
$pop2019Corrected.Range("F5").Value = 5335
$pop2019Corrected.Range("G5").Value = 3004
$pop2019Corrected.Range("H5").Value = 1741
$pop2019Corrected.Range("I5").Value = 566
$pop2019Corrected.Range("J5").Value = 1175

# ---------------------------------------------------------------------
# Sheet: Population 2020-Corrected
# ---------------------------------------------------------------------
# Insert a new row 4 (shifting Arkansas and everything below down by one)
# for Arizona's corrected figures.
$popCorrected2020.Rows.Item(4).Insert()

# Give the new row the same "highlighted / updated" formatting as row 2.
$popCorrected2020.Range("A2:N2").Copy()
$popCorrected2020.Range("A4:N4").PasteSpecial(-4122)

$popCorrected2020.Range("A4").Value = "AZ"
$popCorrected2020.Range("B4").Value = "Arizona"
$popCorrected2020.Range("C4").Value = 37731
$popCorrected2020.Range("D4").Value = 8838
$popCorrected2020.Range("E4").Value = 7897
$popCorrected2020.Range("F4").Value = 5861
$popCorrected2020.Range("G4").Value = 2036
$popCorrected2020.Range("H4").Value = 941
$popCorrected2020.Range("I4").Value = 499
$popCorrected2020.Range("J4").Value = 442
